# edit.ps1 - apply the three textual changes described by the diff:
#   1. "NOME DO PROJETO:" -> "NOME DO PROJETO: STUDIO DE TATUAGEM"
#   2. Append a sentence after the "Iremos priorizar..." paragraph and
#      drop the explicit underline from that paragraph's mark (pilcrow)
#      formatting, matching how Word re-stamps the paragraph-mark rPr
#      after text is typed at the end of the paragraph.
#   3. Merge the split "I" / "nformações..." runs back into a single run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) NOME DO PROJETO: -> NOME DO PROJETO: STUDIO DE TATUAGEM
# ---------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("NOME DO PROJETO:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Collapse(0)
    $r1.InsertAfter(" STUDIO DE TATUAGEM")
}

# ---------------------------------------------------------------------
# 2) Extend the "Iremos priorizar..." paragraph with the extra sentence
# ---------------------------------------------------------------------
$needle2 = "Iremos priorizar a forma de atendimento ao cliente para que ele se sinta familiarizado com nosso estabelecimento e funcionários para que ele passe se sentir confiante com nossos serviços."
$r2 = $d.Content
$found2 = $r2.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Collapse(0)
    $r2.InsertAfter(" Pois nosso studio tende a fazer trabalhos íntimos e pessoais.")
}

# After the insertion, the paragraph mark (pilcrow) that ends this
# paragraph is re-stamped without the explicit single-underline that it
# used to carry (the run text itself stays underlined).
$r2b = $d.Content
$found2b = $r2b.Find.Execute("íntimos e pessoais.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2b) {
    $para = $r2b.Paragraphs.Item(1)
    $markRange = $para.Range
    $markRange.Collapse(0)
    $markRange.Underline = 0
}

# ---------------------------------------------------------------------
# 3) Merge "I" + "nformações que você considera relevante." into one run
# ---------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("Informações que você considera relevante.", $true, $false, $false, $false, $false, $true, 1, $false, "Informações que você considera relevante.", 2)

Write-Output "found1=$found1 found2=$found2 found2b=$found2b found3=$found3"
